$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in B1 (keep it as plain text, not auto-converted to a date serial)
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "2023-02-06"
$ws.Range("B1").Style = "Normal"

# Row 3: only the L-Streak (F3) value changes
$ws.Range("F3").Value = 2

# Row 4: swap winner/loser, new scorer, updated streaks
$ws.Range("B4").Value = "5 Musketeers"
$ws.Range("C4").Value = "Loose Gooses"
$ws.Range("D4").Value = "Sam"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0

# Row 5: swap winner/loser, new scorer, updated streaks
$ws.Range("B5").Value = "5 Musketeers"
$ws.Range("C5").Value = "Wet Willies"
$ws.Range("D5").Value = "Kimmy"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 3

# Remove the old rows 6, 7, 8 (game log now only has 3 entries)
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
